# Atualizei dados da BIBI e da ADD
# Insere dois novos dias (6 e 7) no bloco de Junho/2025, empurrando o restante
# das linhas (Maio, Abril e Março) duas posicoes para baixo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere 2 novas linhas antes da linha 7 (o resto dos dados desliza para baixo)
$ws.Rows("7:8").Insert()

# Novo dia 6 de Junho/2025
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 15533.91
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 2025
$ws.Range("E7").Value = "06/2025"

# Novo dia 7 de Junho/2025
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 19035.4
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 2025
$ws.Range("E8").Value = "06/2025"
